$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# HW9 and HW10 are now graded
$ws.Range("E14").Value = 1.0
$ws.Range("E15").Value = 1.0

# Final Exam total points and "graded?" flag updated now that solutions are posted
$ws.Range("D23").Value = 105.0
$ws.Range("E23").Value = 1.0

# Fill in the Final Exam gradeline column (L) for each letter grade row
$ws.Range("L10").Value = 95.0
$ws.Range("L11").Value = 86.0
$ws.Range("L12").Value = 78.0
$ws.Range("L13").Value = 71.0
$ws.Range("L14").Value = 64.0
$ws.Range("L15").Value = 56.0
$ws.Range("L16").Value = 48.0
$ws.Range("L17").Value = 40.0
$ws.Range("L18").Value = 30.0
$ws.Range("L19").Value = 20.0
$ws.Range("L20").Value = 10.0
$ws.Range("L21").Value = 0.0

$excel.CalculateFullRebuild()
